$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

$row = 13

$ws.Cells.Item($row, 1).Value = "Ik wil mijn bestelling ruilen voor maat M."
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #13: Ik wil mijn bestelling ruilen voor maat M."
$ws.Cells.Item($row, 4).Value = "Retour / Terugbetaling"
$ws.Cells.Item($row, 5).Value = "Beste klant,`nBedankt voor je e-mail. Om je bestelling te ruilen voor maat M, heb ik wat extra informatie nodig. Zou je alsjeblieft de volgende gegevens kunnen doorgeven:`n- Je bestelnummer?`n- Het artikel dat je wilt ruilen en de maat die je wilt ontvangen?`n- Jouw contactgegevens?`nZodra we deze informatie hebben ontvangen, zullen we de ruiling voor je regelen. Mocht je nog andere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam]`nKlantenservice Team`n[Bedrijfsnaam]"
$ws.Cells.Item($row, 6).Value = "2025-07-22 12:39:33"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Nee"
$ws.Cells.Item($row, 9).Value = "Ja"
$ws.Cells.Item($row, 10).Value = "Ja"

# Update the Dashboard summary count for "Retour / Terugbetaling"
$dash.Cells.Item(2, 2).Value = 5

# Extend conditional formatting ranges to cover the new row 13
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = "$($col)2:$($col)12"
    $newRange = "$($col)2:$($col)13"
    $fc = $ws.Range($oldRange).FormatConditions
    if ($fc.Count -gt 0) {
        $fc.Item(1).ModifyAppliesToRange($ws.Range($newRange))
    }
}
